$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '40.743.04'
$ws.Range('E2').Value = '  +3.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.218.91'
$ws.Range('E3').Value = '  +2.26%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.49'
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.634'
$ws.Range('E6').Value = '  +1.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '64.71'
$ws.Range('E7').Value = '  -0.53%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.402'
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0870'
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.545.93'
$ws.Range('E12').Value = '  +2.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '16.00'
$ws.Range('E13').Value = '  -0.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.36'
$ws.Range('E14').Value = '  -0.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.825'
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.61'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.219.21'
$ws.Range('E17').Value = '  +2.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '40.586.05'
$ws.Range('E18').Value = '  +2.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.81'
$ws.Range('E19').Value = '  +1.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0906'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.13'
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '249.88'
$ws.Range('E22').Value = '  +7.44%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.31'
$ws.Range('E25').Value = '  -2.75%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.73'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '173.19'
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.143'
$ws.Range('E28').Value = '  +2.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.48'
$ws.Range('E29').Value = '  +1.83%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.45'
$ws.Range('E30').Value = '  +2.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.87'
$ws.Range('E31').Value = '  +4.71%  '
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.66'
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.77'
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('B35').Value = 'THORChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.10'
$ws.Range('E35').Value = '  -0.72%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.00'
$ws.Range('E36').Value = '  +10.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0632'
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.46'
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.91'
$ws.Range('E40').Value = '  +15.31%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0234'
$ws.Range('E41').Value = '  +0.93%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '104.53'
$ws.Range('E42').Value = '  -0.58%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.45'
$ws.Range('E43').Value = '  +7.60%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.65'
$ws.Range('E44').Value = '  -1.89%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.528.38'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.23'
$ws.Range('E46').Value = '  +1.71%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.12'
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0930'
$ws.Range('E48').Value = '  -0.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.80'
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000198'
$ws.Range('E50').Value = '  +34.07%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '50.41'
$ws.Range('E51').Value = '  +7.31%  '
